$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signup")
$ws.Activate()

# Update the browser value in row 2 from "chrome" to "firefox"
$ws.Range("E2").Value = "firefox"

# Move the active selection to D5
$ws.Range("D5").Select()
